$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Update row 2 values per the diff
$ws.Range("A2").Value = "454-P023454131"
$ws.Range("B2").Value = "CMT"
$ws.Range("C2").Value = "Test 1"

# D2 needs to hold the literal text "true" (not a Boolean). A leading
# apostrophe forces text entry like in the Excel UI; re-apply the Normal
# style afterward so the quote-prefix formatting doesn't linger as an
# explicit cell style.
$ws.Range("D2").Value = "'true"
$ws.Range("D2").Style = "Normal"
